# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Fri Nov 24 02:54:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.460.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.073.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.378.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.074.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.396.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  +4.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0965"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.481.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.19%  "
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.265.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
